$d = $word.ActiveDocument
$d.TrackRevisions = $false

# Remove the existing _GoBack bookmark first so the new one we insert
# (at the paragraph that was last edited) does not collide with it.
try {
    $d.Bookmarks("_GoBack").Delete()
} catch {
}

# Locate the target paragraph by its original text.
$target = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text -like "Cambios f*sicos*") {
        $target = $p
        break
    }
}

$pr = $target.Range
$r = $d.Range($pr.Start, $pr.End - 1)

$xml = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:proofErr w:type="gramStart"/><w:r><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial"/><w:sz w:val="18"/><w:szCs w:val="18"/><w:lang w:val="es-ES_tradnl"/></w:rPr><w:t>cambios</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial"/><w:sz w:val="18"/><w:szCs w:val="18"/><w:lang w:val="es-ES_tradnl"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial"/><w:sz w:val="18"/><w:szCs w:val="18"/><w:lang w:val="es-ES_tradnl"/></w:rPr><w:t>físicos,cambios</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial"/><w:sz w:val="18"/><w:szCs w:val="18"/><w:lang w:val="es-ES_tradnl"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial"/><w:sz w:val="18"/><w:szCs w:val="18"/><w:lang w:val="es-ES_tradnl"/></w:rPr><w:t>químicos,estados</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial"/><w:sz w:val="18"/><w:szCs w:val="18"/><w:lang w:val="es-ES_tradnl"/></w:rPr><w:t xml:space="preserve"> de la </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial"/><w:sz w:val="18"/><w:szCs w:val="18"/><w:lang w:val="es-ES_tradnl"/></w:rPr><w:t>materia,sólidos,líquidos,gases,</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial"/><w:sz w:val="18"/><w:szCs w:val="18"/><w:lang w:val="es-ES_tradnl"/></w:rPr><w:t>temperatura</w:t></w:r><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/><w:proofErr w:type="spellEnd"/></w:p>
'@
$r.InsertXML($xml)
